$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text/percentage/name/link updates (not numeric-looking, Excel keeps them as text) ---
$ws.Range("D2").Value = '64.534.12'
$ws.Range("E2").Value = '  -0.57%  '
$ws.Range("D3").Value = '3.360.38'
$ws.Range("E3").Value = '  -2.38%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("E5").Value = '  -2.60%  '
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("E7").Value = '  -0.53%  '
$ws.Range("D8").Value = '3.354.10'
$ws.Range("E8").Value = '  -2.40%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  +0.91%  '
$ws.Range("E11").Value = '  +2.35%  '
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '3.903.30'
$ws.Range("E15").Value = '  -2.22%  '
$ws.Range("E16").Value = '  +2.00%  '
$ws.Range("E17").Value = '  -1.83%  '
$ws.Range("D18").Value = '3.361.14'
$ws.Range("E18").Value = '  -2.54%  '
$ws.Range("E19").Value = '  -0.05%  '
$ws.Range("D20").Value = '64.494.01'
$ws.Range("E20").Value = '  -0.61%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  +11.99%  '
$ws.Range("E23").Value = '  +10.06%  '
$ws.Range("E24").Value = '  -2.18%  '
$ws.Range("E25").Value = '  +2.34%  '
$ws.Range("E26").Value = '  -1.54%  '
$ws.Range("E27").Value = '  +1.80%  '
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("E29").Value = '  -2.19%  '
$ws.Range("E30").Value = '  +0.22%  '
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("E32").Value = '  -0.40%  '
$ws.Range("E33").Value = '  -0.27%  '
$ws.Range("E34").Value = '  -0.26%  '
$ws.Range("E35").Value = '  -1.92%  '
$ws.Range("E37").Value = '  -7.47%  '
$ws.Range("E38").Value = '  -0.96%  '
$ws.Range("E39").Value = '  -2.73%  '
$ws.Range("D40").Value = '0.0₃0757'
$ws.Range("E40").Value = '  -2.39%  '
$ws.Range("E41").Value = '  +0.56%  '
$ws.Range("D42").Value = '3.104.51'
$ws.Range("E42").Value = '  -2.11%  '
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("E49").Value = '  -1.81%  '
$ws.Range("B44").Value = 'Fetch.AI'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E44").Value = '  +1.82%  '
$ws.Range("B45").Value = 'ThetaToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("E45").Value = '  -4.42%  '
$ws.Range("B46").Value = 'ApeXProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("E46").Value = '  +0.36%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("E47").Value = '  +0.34%  '
$ws.Range("B50").Value = 'Monero'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("E51").Value = '  -0.94%  '

# --- Price (column D) values that look like plain numbers: force text storage so they
#     match the source workbook (values are stored as text, e.g. "1.00", "0.163"),
#     then restore the default "Normal" style so no stray number format sticks around. ---
$numericPriceCells = @("D5","D6","D11","D12","D13","D14","D16","D19","D22","D23","D25","D26","D27","D29","D30","D31","D32","D33","D35","D37","D38","D41","D43","D44","D45","D46","D47","D49","D50","D51")
$numericPriceRange = $ws.Range($numericPriceCells -join ",")
$numericPriceRange.NumberFormat = "@"
$ws.Range("D5").Value = '556.42'
$ws.Range("D6").Value = '175.86'
$ws.Range("D11").Value = '0.163'
$ws.Range("D12").Value = '54.73'
$ws.Range("D13").Value = '0.0000274'
$ws.Range("D14").Value = '9.09'
$ws.Range("D16").Value = '18.45'
$ws.Range("D19").Value = '11.86'
$ws.Range("D22").Value = '456.31'
$ws.Range("D23").Value = '4.90'
$ws.Range("D25").Value = '85.57'
$ws.Range("D26").Value = '13.28'
$ws.Range("D27").Value = '11.00'
$ws.Range("D29").Value = '8.79'
$ws.Range("D30").Value = '29.97'
$ws.Range("D31").Value = '6.61'
$ws.Range("D32").Value = '11.49'
$ws.Range("D33").Value = '583.30'
$ws.Range("D35").Value = '58.67'
$ws.Range("D37").Value = '0.141'
$ws.Range("D38").Value = '35.86'
$ws.Range("D41").Value = '0.378'
$ws.Range("D43").Value = '0.999'
$ws.Range("D44").Value = '2.54'
$ws.Range("D45").Value = '2.80'
$ws.Range("D46").Value = '3.24'
$ws.Range("D47").Value = '0.0411'
$ws.Range("D49").Value = '2.58'
$ws.Range("D50").Value = '137.05'
$ws.Range("D51").Value = '8.37'
$numericPriceRange.Style = "Normal"

